{"js": "// Update the date line and the 25 division problems in the table, in\n// document order. Every run of text in the document changes, so we can\n// walk paragraphs/table cells in order and overwrite each one positionally\n// (this also correctly disambiguates duplicate source texts such as the\n// two \"581\u00f73=\" cells, which map to two different results).\n\nconst newTexts = [\n  \"2025-07-08 Tuesday\",\n  \"484\u00f78=\",\n  \"344\u00f77=\",\n  \"974\u00f77=\",\n  \"549\u00f74=\",\n  \"179\u00f79=\",\n  \"753\u00f75=\",\n  \"306\u00f76=\",\n  \"665\u00f78=\",\n  \"890\u00f75=\",\n  \"753\u00f78=\",\n  \"678\u00f77=\",\n  \"585\u00f76=\",\n  \"462\u00f74=\",\n  \"784\u00f74=\",\n  \"834\u00f73=\",\n  \"135\u00f78=\",\n  \"558\u00f78=\",\n  \"684\u00f72=\",\n  \"842\u00f78=\",\n  \"974\u00f74=\",\n  \"328\u00f78=\",\n  \"291\u00f79=\",\n  \"787\u00f72=\",\n  \"988\u00f74=\",\n  \"428\u00f78=\",\n];\n\nlet idx = 0;\n\n// 1) Title paragraph (first paragraph of the body, outside of the table).\nconst body = context.document.body;\nconst titleParas = body.paragraphs;\ntitleParas.load(\"items\");\nawait context.sync();\n\ntitleParas.items[0].insertText(newTexts[idx], Word.InsertLocation.replace);\nidx++;\n\n// 2) Table cells, row by row, left to right \u2014 only the rows that actually\n// contain text (the diff never touches the intentionally blank rows).\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.rows.load(\"items\");\nawait context.sync();\n\nfor (const row of table.rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nfor (const row of table.rows.items) {\n  for (const cell of row.cells.items) {\n    cell.body.paragraphs.load(\"items\");\n  }\n}\nawait context.sync();\n\nfor (const row of table.rows.items) {\n  for (const cell of row.cells.items) {\n    const paras = cell.body.paragraphs.items;\n    if (paras.length === 0) continue;\n    // Blank rows have an empty paragraph with no runs/text \u2014 skip them.\n    const text = paras[0].text;\n    if (text === \"\") continue;\n\n    paras[0].insertText(newTexts[idx], Word.InsertLocation.replace);\n    idx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 25 division problems in the table, in\n# document order. Every text run in the document changes, so each cell is\n# addressed positionally by (row, column) \u2014 this correctly disambiguates\n# duplicate source texts such as the two \"581\u00f73=\" cells, which map to two\n# different results (\"135\u00f78=\" and \"328\u00f78=\").\n\n$d = $word.ActiveDocument\n\n# 1) Title paragraph (first paragraph of the body, outside of the table).\n$d.Paragraphs.Item(1).Range.Text = \"2025-07-08 Tuesday\"\n\n# 2) Table cells, row by row, left to right \u2014 only the rows that actually\n# contain text (the diff never touches the intentionally blank rows).\n$newTexts = @(\n  \"484\u00f78=\", \"344\u00f77=\", \"974\u00f77=\", \"549\u00f74=\", \"179\u00f79=\",\n  \"753\u00f75=\", \"306\u00f76=\", \"665\u00f78=\", \"890\u00f75=\", \"753\u00f78=\",\n  \"678\u00f77=\", \"585\u00f76=\", \"462\u00f74=\", \"784\u00f74=\", \"834\u00f73=\",\n  \"135\u00f78=\", \"558\u00f78=\", \"684\u00f72=\", \"842\u00f78=\", \"974\u00f74=\",\n  \"328\u00f78=\", \"291\u00f79=\", \"787\u00f72=\", \"988\u00f74=\", \"428\u00f78=\"\n)\n\n$t = $d.Tables.Item(1)\n$idx = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n  for ($c = 1; $c -le $t.Columns.Count; $c++) {\n    $cell = $t.Cell($r, $c)\n    $txt = $cell.Range.Text\n    # $txt always carries a trailing cell-mark (length 2 when the cell is\n    # otherwise empty), so compare length rather than equality with \"\".\n    if ($txt.Length -gt 2) {\n      $cell.Range.Text = $newTexts[$idx]\n      $idx++\n    }\n  }\n}\n"}
